$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A3").Value = "2026-01-12 16:18:11"
$ws.Range("B3").Value = 200
$ws.Range("C3").Value = "food"
